$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.785.60"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "2.353.16"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.08%  "
$c = $ws.Range("D5")
$c.Value = "'239.95"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$c = $ws.Range("D6")
$c.Value = "'0.669"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.62%  "
$c = $ws.Range("D7")
$c.Value = "'73.25"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$c = $ws.Range("D9")
$c.Value = "'0.602"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  +1.66%  "
$c = $ws.Range("D11")
$c.Value = "'60.78"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'35.62"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +10.04%  "
$ws.Range("E13").Value = "  +0.24%  "
$c = $ws.Range("D14")
$c.Value = "'7.20"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.83%  "
$c = $ws.Range("D15")
$c.Value = "'16.22"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.35%  "
$c = $ws.Range("D16")
$c.Value = "'0.910"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "2.360.75"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "43.738.71"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  +1.43%  "
$c = $ws.Range("D20")
$c.Value = "'77.74"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.99%  "
$c = $ws.Range("D21")
$c.Value = "'6.54"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.08%  "
$c = $ws.Range("D22")
$c.Value = "'252.63"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("E23").Value = "  +0.04%  "
$c = $ws.Range("D24")
$c.Value = "'3.77"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("E26").Value = "  -0.32%  "
$c = $ws.Range("D27")
$c.Value = "'10.49"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.93%  "
$c = $ws.Range("D28")
$c.Value = "'2.29"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.06%  "
$c = $ws.Range("D29")
$c.Value = "'175.59"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$c = $ws.Range("D30")
$c.Value = "'22.27"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  -2.21%  "
$c = $ws.Range("D33")
$c.Value = "'0.0747"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("E34").Value = "  -3.70%  "
$c = $ws.Range("D35")
$c.Value = "'5.34"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.23%  "
$c = $ws.Range("D36")
$c.Value = "'3.77"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "
$c = $ws.Range("D37")
$c.Value = "'6.60"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.70%  "
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("E39").Value = "  -1.48%  "
$c = $ws.Range("D40")
$c.Value = "'5.48"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +13.60%  "
$c = $ws.Range("D41")
$c.Value = "'65.13"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +12.07%  "
$c = $ws.Range("D42")
$c.Value = "'19.81"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D44")
$c.Value = "'9.02"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D45")
$c.Value = "'0.106"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.12%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D46")
$c.Value = "'2.48"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D47")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  -1.92%  "
$c = $ws.Range("D50")
$c.Value = "'97.80"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("E51").Value = "  +2.03%  "
